$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.820.90"
$ws.Range("E2").Value = "  -1.31%  "
$ws.Range("D3").Value = "3.571.42"
$ws.Range("E3").Value = "  -1.66%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'577.17"
$ws.Range("E5").Value = "  -2.87%  "
$ws.Range("D6").Value = "'188.57"
$ws.Range("E6").Value = "  -2.01%  "
$ws.Range("E7").Value = "  -3.93%  "
$ws.Range("D8").Value = "3.570.61"
$ws.Range("E8").Value = "  -1.28%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "'0.177"
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("D11").Value = "'0.660"
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("D12").Value = "'55.84"
$ws.Range("E12").Value = "  -4.02%  "
$ws.Range("E13").Value = "  +2.19%  "
$ws.Range("D14").Value = "'9.58"
$ws.Range("E14").Value = "  -2.16%  "
$ws.Range("D15").Value = "4.147.90"
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("D16").Value = "'19.72"
$ws.Range("D17").Value = "3.564.75"
$ws.Range("E17").Value = "  -1.72%  "
$ws.Range("D18").Value = "69.767.69"
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("D19").Value = "'12.60"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("D22").Value = "'475.73"
$ws.Range("E22").Value = "  -4.01%  "
$ws.Range("D23").Value = "'19.21"
$ws.Range("E23").Value = "  +11.84%  "
$ws.Range("E24").Value = "  -7.78%  "
$ws.Range("D25").Value = "'95.76"
$ws.Range("E25").Value = "  +4.94%  "
$ws.Range("D26").Value = "'4.37"
$ws.Range("E26").Value = "  -2.82%  "
$ws.Range("D27").Value = "'3.01"
$ws.Range("E27").Value = "  -3.86%  "
$ws.Range("D28").Value = "'10.97"
$ws.Range("E28").Value = "  -2.94%  "
$ws.Range("D29").Value = "'9.31"
$ws.Range("E29").Value = "  -1.86%  "
$ws.Range("D30").Value = "'32.29"
$ws.Range("E30").Value = "  -0.33%  "
$ws.Range("D31").Value = "'7.72"
$ws.Range("E31").Value = "  +1.18%  "
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("D33").Value = "'12.15"
$ws.Range("E33").Value = "  -1.10%  "
$ws.Range("D34").Value = "'66.12"
$ws.Range("E34").Value = "  +1.26%  "
$ws.Range("D35").Value = "'584.53"
$ws.Range("E35").Value = "  -5.52%  "
$ws.Range("D36").Value = "'38.93"
$ws.Range("E36").Value = "  +1.83%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").Value = "0.0₃0794"
$ws.Range("E38").Value = "  -4.41%  "
$ws.Range("D39").Value = "'0.394"
$ws.Range("E39").Value = "  -3.42%  "
$ws.Range("D40").Value = "'3.20"
$ws.Range("E40").Value = "  +15.97%  "
$ws.Range("E41").Value = "  -7.00%  "
$ws.Range("D42").Value = "'3.45"
$ws.Range("E42").Value = "  -5.32%  "
$ws.Range("D43").Value = "3.221.10"
$ws.Range("E43").Value = "  -3.64%  "
$ws.Range("D44").Value = "'2.84"
$ws.Range("E44").Value = "  +4.88%  "
$ws.Range("D45").Value = "'3.08"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").Value = "'0.0440"
$ws.Range("E46").Value = "  -2.12%  "
$ws.Range("D47").Value = "'3.39"
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("D48").Value = "'9.47"
$ws.Range("E48").Value = "  +2.50%  "
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("D50").Value = "'0.998"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").Value = "'3.13"
$ws.Range("E51").Value = "  -7.05%  "
